$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added at the top of the date-ordered list for this
# market/product, pushing every existing data row (14-26) down by one
# (they become rows 15-27). Inserting a row at 14 reproduces that shift.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new record's values.
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(14, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(14, 4).Value = 44762
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100112013
$ws.Cells.Item(14, 7).Value = 'Alcachofa'
$ws.Cells.Item(14, 8).Value = 'Madrigal'
$ws.Cells.Item(14, 9).Value = 'Primera'
$ws.Cells.Item(14, 10).Value = 200
$ws.Cells.Item(14, 11).Value = 19000
$ws.Cells.Item(14, 12).Value = 20000
$ws.Cells.Item(14, 13).Value = 19500
$ws.Cells.Item(14, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(14, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(14, 16).Value = 488
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(14, 18).Value = 'Hortaliza'
